$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rename the sheet
$ws.Name = "ESTADÍSITICA RES"

# Add the defined name TOTAL_PRUEBAS -> D4:D28 on the renamed sheet
$wb.Names.Add("TOTAL_PRUEBAS", "='ESTADÍSITICA RES'!`$D`$4:`$D`$28")

# Update a couple of the descriptive labels in column H
$ws.Range("H16").Value = "La posición ocupada a partir de la nota Total Pruebas"
$ws.Range("H20").Value = "La nota que aparece con más frecuencia en cada prueba"

# The TOTAL PRUEBAS column (D4:D28) no longer holds the ROUND formula;
# clear it out and mark it with the "Output" (Salida) cell style instead.
$ws.Range("D4:D28").ClearContents()
$ws.Range("D4:D28").Style = "Output"

# Move the active selection as recorded in the workbook
[void]$ws.Activate()
$ws.Range("F35").Select() | Out-Null

Write-Host "Edits applied"
